$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Version (row 3) and Date (row 8) values in the Metadata sheet.
$ws.Cells.Item(3, 2).Value = "2.0.1-sd-202510-matchbox-patch"
$ws.Cells.Item(8, 2).Value = "2025-10-29T22:15:57+01:00"

# Insert a new "Jurisdiction" property row after "Contact" (row 10), before
# "Description" (previously row 11), pushing all following rows down by one.
$ws.Rows.Item(11).Insert()

# Copy the formatting (style) from the row above so the new row matches the
# rest of the table instead of getting a brand-new default style.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""
